# feat: patch assign simple
#
# 1) Update the scan metadata (target port + start date).
# 2) Clear the stale placeholder "Endpoint" values (the literal "*",
#    "POST http://192.168.1.14:8080/", "GET http://192.168.1.14:8080/"
#    filler text) out of column D for every data row, leaving the column's
#    formatting/style untouched. Header rows (whose column D reads
#    "Endpoint") are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- metadata updates -------------------------------------------------
$ws.Range("A3").Value = "Target URL/IP : 192.168.1.14:3000"
$ws.Range("C3").Value = "Start Date : 2023-06-12"

# --- clear stale Endpoint placeholder values in column D --------------
$clearRanges = @(
    "D7:D16",
    "D19:D31",
    "D35:D37",
    "D41:D43",
    "D47",
    "D49:D50",
    "D61:D64",
    "D68",
    "D72:D77",
    "D80",
    "D84:D89",
    "D95:D96",
    "D99:D100",
    "D102",
    "D105:D111",
    "D117",
    "D119",
    "D125"
)

foreach ($addr in $clearRanges) {
    $ws.Range($addr).ClearContents()
}
